$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-safe updates (values that Excel will not auto-convert to a number) ---
$ws.Range('D2').Value = '62.416.52'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '2.453.34'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('E5').Value = '  +3.33%  '
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.38%  '
$ws.Range('D9').Value = '2.449.19'
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('E10').Value = '  +1.12%  '
$ws.Range('E11').Value = '  +2.38%  '
$ws.Range('E12').Value = '  +0.79%  '
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').Value = '2.895.96'
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('D17').Value = '62.215.38'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '2.454.92'
$ws.Range('E19').Value = '  -2.77%  '
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('E21').Value = '  +0.91%  '
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('E23').Value = '  -3.04%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('E26').Value = '  +2.11%  '
$ws.Range('E27').Value = '  -5.31%  '
$ws.Range('D28').Value = '0.0₃0976'
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('D29').Value = '2.574.77'
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('E31').Value = '  -1.76%  '
$ws.Range('E32').Value = '  -0.73%  '
$ws.Range('E33').Value = '  +1.68%  '
$ws.Range('E34').Value = '  +0.87%  '
$ws.Range('E35').Value = '  -1.77%  '
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('E38').Value = '  +0.82%  '
$ws.Range('E39').Value = '  +5.36%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('E40').Value = '  +1.09%  '
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E41').Value = '  -1.59%  '
$ws.Range('E42').Value = '  +1.47%  '
$ws.Range('E43').Value = '  -1.15%  '
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('D48').Value = '0.0₆0255'
$ws.Range('E48').Value = '  +14.81%  '
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('E51').Value = '  -1.93%  '

# --- Numeric-looking text updates: force text storage so they remain strings, not numbers ---
$forceTextCells = @('D4', 'D5', 'D6', 'D14', 'D15', 'D20', 'D21', 'D25', 'D26', 'D27', 'D32', 'D38', 'D39', 'D40', 'D41', 'D42', 'D45', 'D46', 'D47', 'D51')
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D4').Value = '0.998'
$ws.Range('D5').Value = '583.80'
$ws.Range('D6').Value = '143.57'
$ws.Range('D14').Value = '26.52'
$ws.Range('D15').Value = '0.0000178'
$ws.Range('D20').Value = '7.15'
$ws.Range('D21').Value = '327.68'
$ws.Range('D25').Value = '65.79'
$ws.Range('D26').Value = '9.20'
$ws.Range('D27').Value = '591.13'
$ws.Range('D32').Value = '8.01'
$ws.Range('D38').Value = '0.379'
$ws.Range('D39').Value = '153.86'
$ws.Range('D40').Value = '5.31'
$ws.Range('D41').Value = '18.44'
$ws.Range('D42').Value = '42.89'
$ws.Range('D45').Value = '2.50'
$ws.Range('D46').Value = '142.99'
$ws.Range('D47').Value = '3.66'
$ws.Range('D51').Value = '19.89'
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}
